$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.304.50"
$ws.Range("E2").Value = "  +0.68%  "
$ws.Range("D3").Value = "3.616.15"
$ws.Range("E3").Value = "  +2.31%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'601.98"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "'195.70"
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.212"
$ws.Range("E9").Value = "  +4.09%  "
$ws.Range("D10").Value = "'0.647"
$ws.Range("E10").Value = "  -0.54%  "
$ws.Range("D11").Value = "'53.27"
$ws.Range("E11").Value = "  -0.84%  "
$ws.Range("E12").Value = "  +0.75%  "
$ws.Range("E13").Value = "  +0.55%  "
$ws.Range("D14").Value = "4.186.57"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").Value = "'601.55"
$ws.Range("E15").Value = "  -0.47%  "
$ws.Range("D16").Value = "'12.98"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("D17").Value = "70.406.18"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "3.615.92"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("E19").Value = "  -0.34%  "
$ws.Range("E20").Value = "  +1.70%  "
$ws.Range("E21").Value = "  +0.94%  "
$ws.Range("D22").Value = "'18.74"
$ws.Range("E22").Value = "  +3.03%  "
$ws.Range("E23").Value = "  -1.40%  "
$ws.Range("D24").Value = "'103.16"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -4.60%  "
$ws.Range("E27").Value = "  -2.71%  "
$ws.Range("D28").Value = "'9.74"
$ws.Range("E28").Value = "  +1.14%  "
$ws.Range("D29").Value = "'33.84"
$ws.Range("E29").Value = "  +0.93%  "
$ws.Range("D30").Value = "'4.72"
$ws.Range("E30").Value = "  +8.98%  "
$ws.Range("D31").Value = "'7.30"
$ws.Range("E31").Value = "  +2.64%  "
$ws.Range("E32").Value = "  -2.19%  "
$ws.Range("E33").Value = "  +2.41%  "
$ws.Range("D34").Value = "'63.35"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Value = "0.0₃0887"
$ws.Range("E35").Value = "  +3.31%  "
$ws.Range("D36").Value = "3.937.82"
$ws.Range("E36").Value = "  +5.08%  "
$ws.Range("D37").Value = "'533.67"
$ws.Range("E37").Value = "  +9.73%  "
$ws.Range("E38").Value = "  +0.15%  "
$ws.Range("D39").Value = "'3.04"
$ws.Range("E39").Value = "  +0.47%  "
$ws.Range("D40").Value = "'36.86"
$ws.Range("E40").Value = "  +0.77%  "
$ws.Range("E41").Value = "  -0.81%  "
$ws.Range("E42").Value = "  -2.09%  "
$ws.Range("E43").Value = "  +0.55%  "
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("D45").Value = "'3.57"
$ws.Range("E45").Value = "  +8.26%  "
$ws.Range("E46").Value = "  +1.31%  "
$ws.Range("D47").Value = "'0.140"
$ws.Range("E47").Value = "  -0.33%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("E49").Value = "  -0.24%  "
$ws.Range("D50").Value = "'0.000249"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "'1.30"
$ws.Range("E51").Value = "  +1.35%  "
